$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.300.06"
$ws.Range("E2").Value = "  +1.40%  "
$ws.Range("D3").Value = "3.913.78"
$ws.Range("E3").Value = "  -0.73%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "488.84"
$ws.Range("E5").Value = "  +3.90%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "146.89"
$ws.Range("E6").Value = "  +0.52%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.622"
$ws.Range("E7").Value = "  -0.41%  "
$ws.Range("E8").Value = "  -0.14%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.728"
$ws.Range("E9").Value = "  -0.45%  "
$ws.Range("E10").Value = "  -3.66%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0000346"
$ws.Range("E11").Value = "  -3.07%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "43.04"
$ws.Range("E12").Value = "  -0.59%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "10.75"
$ws.Range("E13").Value = "  +3.67%  "
$ws.Range("D14").Value = "4.531.78"
$ws.Range("E14").Value = "  -0.94%  "
$ws.Range("D15").Value = "3.912.85"
$ws.Range("E15").Value = "  -0.62%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.19"
$ws.Range("E16").Value = "  -6.22%  "
$ws.Range("E17").Value = "  -1.10%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "19.86"
$ws.Range("E18").Value = "  +0.03%  "
$ws.Range("E19").Value = "  -1.69%  "
$ws.Range("D20").Value = "68.327.36"
$ws.Range("E20").Value = "  +1.04%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "433.84"
$ws.Range("E21").Value = "  -0.20%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.53"
$ws.Range("E22").Value = "  +4.38%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "14.93"
$ws.Range("E23").Value = "  +3.65%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "87.70"
$ws.Range("E24").Value = "  +0.28%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "11.45"
$ws.Range("E25").Value = "  +16.18%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.28"
$ws.Range("E26").Value = "  +10.84%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "3.61"
$ws.Range("E27").Value = "  -0.33%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "38.20"
$ws.Range("E28").Value = "  -1.02%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.71"
$ws.Range("E29").Value = "  +0.53%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "726.05"
$ws.Range("E30").Value = "  +0.67%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "13.76"
$ws.Range("E31").Value = "  +2.57%  "
$ws.Range("E32").Value = "  -1.35%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.91"
$ws.Range("E33").Value = "  +3.74%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.27"
$ws.Range("E34").Value = "  +17.86%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "41.63"
$ws.Range("E35").Value = "  -1.31%  "
$ws.Range("D36").Value = "0.0₃0867"
$ws.Range("E36").Value = "  +1.68%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "60.61"
$ws.Range("E37").Value = "  +4.62%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.405"
$ws.Range("E38").Value = "  +20.80%  "
$ws.Range("E39").Value = "  -1.79%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.999"
$ws.Range("E40").Value = "  +0.01%  "
$ws.Range("E41").Value = "  +17.14%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0481"
$ws.Range("E42").Value = "  +1.21%  "
$ws.Range("E43").Value = "  +3.51%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.92"
$ws.Range("E44").Value = "  +3.39%  "
$ws.Range("E45").Value = "  -0.81%  "
$ws.Range("E46").Value = "  -0.04%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.32"
$ws.Range("E47").Value = "  +5.23%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.41"
$ws.Range("E48").Value = "  -3.16%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.13"
$ws.Range("E49").Value = "  -3.36%  "
$ws.Range("D50").Value = "0.0₆0345"
$ws.Range("E50").Value = "  +33.50%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "144.68"
$ws.Range("E51").Value = "  -2.18%  "
